$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '59.009.05'
Set-TextValue $ws.Range('E2') '  +4.61%  '
Set-TextValue $ws.Range('D3') '3.313.39'
Set-TextValue $ws.Range('E3') '  +2.09%  '
Set-TextValue $ws.Range('E4') '  +0.07%  '
Set-TextValue $ws.Range('D5') '408.09'
Set-TextValue $ws.Range('E5') '  +2.44%  '
Set-TextValue $ws.Range('D6') '110.60'
Set-TextValue $ws.Range('E6') '  -0.66%  '
Set-TextValue $ws.Range('E7') '  +4.68%  '
Set-TextValue $ws.Range('D8') '0.999'
Set-TextValue $ws.Range('E8') '  +0.01%  '
Set-TextValue $ws.Range('D9') '0.633'
Set-TextValue $ws.Range('E9') '  +1.80%  '
Set-TextValue $ws.Range('D10') '39.72'
Set-TextValue $ws.Range('E10') '  +0.97%  '
Set-TextValue $ws.Range('D11') '0.0979'
Set-TextValue $ws.Range('E11') '  +4.18%  '
Set-TextValue $ws.Range('E12') '  +1.17%  '
Set-TextValue $ws.Range('D13') '3.840.95'
Set-TextValue $ws.Range('E13') '  +2.71%  '
Set-TextValue $ws.Range('D14') '8.40'
Set-TextValue $ws.Range('E14') '  +3.43%  '
Set-TextValue $ws.Range('D15') '19.33'
Set-TextValue $ws.Range('E15') '  +0.79%  '
Set-TextValue $ws.Range('D16') '3.308.59'
Set-TextValue $ws.Range('E16') '  +2.42%  '
Set-TextValue $ws.Range('E17') '  -0.49%  '
Set-TextValue $ws.Range('D18') '58.890.00'
Set-TextValue $ws.Range('E18') '  +4.69%  '
Set-TextValue $ws.Range('D19') '10.70'
Set-TextValue $ws.Range('E19') '  -3.28%  '
Set-TextValue $ws.Range('D20') '3.31'
Set-TextValue $ws.Range('E20') '  -1.16%  '
Set-TextValue $ws.Range('E21') '  +4.27%  '
Set-TextValue $ws.Range('D22') '12.94'
Set-TextValue $ws.Range('E22') '  -1.60%  '
Set-TextValue $ws.Range('D23') '302.79'
Set-TextValue $ws.Range('E23') '  +1.56%  '
Set-TextValue $ws.Range('D24') '74.86'
Set-TextValue $ws.Range('E24') '  -1.27%  '
Set-TextValue $ws.Range('D25') '3.22'
Set-TextValue $ws.Range('E25') '  -0.48%  '
Set-TextValue $ws.Range('D26') '28.49'
Set-TextValue $ws.Range('E26') '  +1.01%  '
Set-TextValue $ws.Range('D27') '4.46'
Set-TextValue $ws.Range('E27') '  +2.30%  '
Set-TextValue $ws.Range('D28') '7.83'
Set-TextValue $ws.Range('E28') '  -4.62%  '
Set-TextValue $ws.Range('E29') '  -0.49%  '
Set-TextValue $ws.Range('D30') '7.30'
Set-TextValue $ws.Range('E30') '  -1.03%  '
Set-TextValue $ws.Range('E31') '  -0.03%  '
Set-TextValue $ws.Range('E32') '  +1.33%  '
Set-TextValue $ws.Range('D33') '11.39'
Set-TextValue $ws.Range('E33') '  +1.89%  '
Set-TextValue $ws.Range('D34') '40.25'
Set-TextValue $ws.Range('E34') '  +8.71%  '
Set-TextValue $ws.Range('D35') '0.0528'
Set-TextValue $ws.Range('E35') '  +7.40%  '
Set-TextValue $ws.Range('E36') '  -0.06%  '
Set-TextValue $ws.Range('D37') '51.86'
Set-TextValue $ws.Range('E37') '  +0.85%  '
Set-TextValue $ws.Range('D38') '3.23'
Set-TextValue $ws.Range('E38') '  +4.24%  '
Set-TextValue $ws.Range('E40') '  -2.57%  '
Set-TextValue $ws.Range('D41') '137.93'
Set-TextValue $ws.Range('E41') '  +2.17%  '
Set-TextValue $ws.Range('E42') '  +1.58%  '
Set-TextValue $ws.Range('E44') '  -2.44%  '
Set-TextValue $ws.Range('D45') '16.66'
Set-TextValue $ws.Range('E45') '  -5.34%  '
Set-TextValue $ws.Range('E46') '  -2.62%  '
Set-TextValue $ws.Range('D47') '2.27'
Set-TextValue $ws.Range('E47') '  +8.24%  '
Set-TextValue $ws.Range('D48') '22.24'
Set-TextValue $ws.Range('E48') '  -0.79%  '
Set-TextValue $ws.Range('D49') '2.181.81'
Set-TextValue $ws.Range('E49') '  +1.95%  '
Set-TextValue $ws.Range('E50') '  +0.01%  '
Set-TextValue $ws.Range('B51') 'THORChain'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range('D51') '6.35'
Set-TextValue $ws.Range('E51') '  +5.54%  '
